# Apply updated crypto price/volume data to Sheet1 (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    # Prefix with an apostrophe so Excel stores the value as literal text
    # (matches the workbook's existing inlineStr/text cells) instead of
    # auto-coercing numeric-looking strings (e.g. '294.14') into numbers.
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# Row 2
Set-TextCell "D2" '43.899.79'
Set-TextCell "E2" '  -1.16%  '

# Row 3
Set-TextCell "D3" '2.192.86'
Set-TextCell "E3" '  -2.30%  '

# Row 4
Set-TextCell "E4" '  -0.08%  '

# Row 5
Set-TextCell "D5" '294.14'
Set-TextCell "E5" '  -4.05%  '

# Row 6
Set-TextCell "D6" '87.70'
Set-TextCell "E6" '  -5.90%  '

# Row 7
Set-TextCell "E7" '  -1.32%  '

# Row 8
Set-TextCell "E8" '  -0.09%  '

# Row 9
Set-TextCell "D9" '0.476'
Set-TextCell "E9" '  -8.94%  '

# Row 10
Set-TextCell "D10" '32.02'
Set-TextCell "E10" '  -7.37%  '

# Row 11
Set-TextCell "D11" '0.0763'
Set-TextCell "E11" '  -5.83%  '

# Row 12
Set-TextCell "E12" '  -1.69%  '

# Row 13
Set-TextCell "D13" '6.71'
Set-TextCell "E13" '  -6.07%  '

# Row 14
Set-TextCell "D14" '2.528.03'
Set-TextCell "E14" '  -2.25%  '

# Row 15
Set-TextCell "D15" '2.251.48'
Set-TextCell "E15" '  -6.00%  '

# Row 16
Set-TextCell "E16" '  -5.31%  '

# Row 17
Set-TextCell "D17" '0.761'
Set-TextCell "E17" '  -9.01%  '

# Row 18
Set-TextCell "D18" '43.470.91'
Set-TextCell "E18" '  -1.38%  '

# Row 19
Set-TextCell "D19" '0.0₃0877'
Set-TextCell "E19" '  -8.81%  '

# Row 20
Set-TextCell "D20" '5.78'
Set-TextCell "E20" '  -9.26%  '

# Row 21
Set-TextCell "D21" '10.64'
Set-TextCell "E21" '  -14.10%  '

# Row 22
Set-TextCell "D22" '62.46'
Set-TextCell "E22" '  -4.93%  '

# Row 23
Set-TextCell "D23" '229.35'
Set-TextCell "E23" '  -3.22%  '

# Row 24
Set-TextCell "E24" '  -11.88%  '

# Row 25
Set-TextCell "E25" '  -0.04%  '

# Row 26
Set-TextCell "E26" '  -8.86%  '

# Row 27
Set-TextCell "E27" '  +0.30%  '

# Row 28
Set-TextCell "D28" '9.14'
Set-TextCell "E28" '  -6.93%  '

# Row 29
Set-TextCell "D29" '35.09'
Set-TextCell "E29" '  -9.43%  '

# Row 30
Set-TextCell "D30" '19.04'
Set-TextCell "E30" '  -5.04%  '

# Row 31
Set-TextCell "D31" '145.68'
Set-TextCell "E31" '  -5.20%  '

# Row 33
Set-TextCell "E33" '  -5.44%  '

# Row 34
Set-TextCell "D34" '0.0722'
Set-TextCell "E34" '  -9.44%  '

# Row 35
Set-TextCell "E35" '  -3.65%  '

# Row 36
Set-TextCell "D36" '2.84'
Set-TextCell "E36" '  -8.65%  '

# Row 37
Set-TextCell "E37" '  -6.71%  '

# Row 38
Set-TextCell "D38" '1.63'
Set-TextCell "E38" '  -10.40%  '

# Row 39
Set-TextCell "B39" 'RenderToken'
Set-TextCell "C39" 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell "D39" '3.47'
Set-TextCell "E39" '  -9.07%  '

# Row 40
Set-TextCell "E40" '  -7.75%  '

# Row 41
Set-TextCell "B41" 'Celestia'
Set-TextCell "C41" 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextCell "D41" '13.10'
Set-TextCell "E41" '  -9.95%  '

# Row 42
Set-TextCell "D42" '1.01'
Set-TextCell "E42" '  -0.26%  '

# Row 43
Set-TextCell "D43" '3.02'
Set-TextCell "E43" '  -12.21%  '

# Row 44
Set-TextCell "D44" '1.749.05'
Set-TextCell "E44" '  +0.73%  '

# Row 45
Set-TextCell "D45" '1.62'
Set-TextCell "E45" '  +1.82%  '

# Row 46
Set-TextCell "D46" '14.18'
Set-TextCell "E46" '  -0.50%  '

# Row 47
Set-TextCell "D47" '72.15'
Set-TextCell "E47" '  -10.07%  '

# Row 48
Set-TextCell "D48" '65.58'
Set-TextCell "E48" '  -5.78%  '

# Row 49
Set-TextCell "B49" 'HuobiToken'
Set-TextCell "C49" 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextCell "D49" '2.74'
Set-TextCell "E49" '  +8.93%  '

# Row 50
Set-TextCell "B50" 'Algorand'
Set-TextCell "C50" 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextCell "D50" '0.170'
Set-TextCell "E50" '  -11.45%  '

# Row 51
Set-TextCell "D51" '90.84'
Set-TextCell "E51" '  -8.73%  '
